$d = $word.ActiveDocument

# Namespace prefix used for the raw OOXML fragments we insert below.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Build the block of new paragraphs that goes after the existing last
# paragraph (the one ending "...CAS."). The first new paragraph is a
# completely empty one (no run at all), followed by the date and the
# four new diary entries. The bookmark is intentionally left out of the
# raw fragment here - it gets re-added (and the stale one at the old
# location automatically removed) by Bookmarks.Add further down.
$newParasXml =
    "<w:p $wns/>" +
    "<w:p $wns><w:r><w:t>27/01/14</w:t></w:r></w:p>" +
    "<w:p $wns><w:r><w:t>Researched several papers on the topic of concurrent data structures, placed them into the research folder</w:t></w:r></w:p>" +
    "<w:p $wns><w:r><w:t>Amended the output of the program to be easier to paste into excel for graph generation</w:t></w:r></w:p>" +
    "<w:p $wns><w:r><w:t>Gathered data from stoker for both locked and spinlock modes of ring buffer</w:t></w:r></w:p>" +
    "<w:p $wns><w:r><w:t>Need to gather data from my machine, spoon and ducss</w:t></w:r></w:p>"

# Insert at a collapsed range sitting at the very end of the document
# body (right after the existing "_GoBack" bookmark on the CAS. line).
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)
$insertionPoint.InsertXML($newParasXml)

# The "_GoBack" bookmark used to sit around the end of the "...CAS."
# paragraph; re-home it onto the newly added "Researched..." paragraph
# (excluding its trailing paragraph mark), matching where Word leaves
# it after the most recent block of typed text. Adding a bookmark with
# a name that already exists moves it (the old one is removed).
$target = $d.Paragraphs.Item(29)
$bookmarkRange = $d.Range($target.Range.Start, $target.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
